$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.10186050314992201,
    -0.0059999999472388765,
    -0.0039999999570348166,
    -0.0079999999180273562,
    -0.0029999999570771152,
    -0.0019999999554816128,
    -0.0099999998850148764,
    -0.0099999998808000257,
    -0.0019999999450477368,
    -0.0019999999406365987,
    -0.0029999999318883752,
    -0.0034999999265550308,
    -0.003499999922434327,
    -0.0079999998825108776,
    -0.00099999994243393786,
    -0.0019999999334188168,
    -0.0019999999331039575,
    0.00042755312226283593,
    -0.0039999999644799722,
    0.014291144432881708,
    -0.0039999999647601925,
    -0.0039999999645052853,
    -0.0049999999479171109,
    0.01733457757450374,
    -0.019999999815699887,
    -0.002499999950629217,
    0.0024109378181975671,
    -0.0019999999505033728,
    -0.006999999904958365,
    0.010935905929813305,
    -0.0069999999005858626,
    -0.0099999998740507579,
    -0.0039999999257176455
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
